$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column H1, copying the style/format of the existing
# header cells (e.g. G1) so it matches the bold/bordered header formatting.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the new data values in H2 and H3 (plain numeric cells, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
